$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# Row 16: Queue - LC75+ entry
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 5).Value = "933, 649"
$ws.Cells.Item(16, 6).Value = "Med"

# Row 17: Linked List entry - label tweak (drop the trailing '+') and quantity fix
$ws.Cells.Item(17, 2).Value = "Linked List - LC75"
$ws.Cells.Item(17, 3).Value = 1

# Row 18: Linked List entry - quantity fix
$ws.Cells.Item(18, 3).Value = 1

# Move the active selection to F24, matching where work left off
$ws.Range("F24").Select()
